$d = $word.ActiveDocument

# The last paragraph in the document is the existing limitation that ends with
# "Just tested with Firefox version 30.0 and 31.0, not sure older Firefox's
# version can run it." Append a brand-new list paragraph right after it (same
# ListParagraph style / numbering) containing the new limitation text.

$target = $d.Paragraphs.Last
$r = $target.Range
$r.Collapse(0)            # wdCollapseEnd - move to the very end of that paragraph
$r.InsertParagraphAfter() # creates a new (empty) paragraph that inherits the style/numbering
$r.InsertAfter("Does not check API Key before scanning.")
$r.Collapse(0)            # collapse to just after the newly-typed text

# Word keeps the "_GoBack" bookmark anchored at the location of the most
# recent edit, so move it from the old last paragraph onto the new one.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $r)
